# Clarify service names and usage in the architecture diagram on slide 2.
#
# 1. Rename the "ProductCustomer" cloud to "Ownership Service".
# 2. Reverse/clarify the arrowheads on the two connectors that run from
#    "Atlas Reader" down into the two cloud shapes, switching them from a
#    tail-only arrow to a head-side triangle (the second one using the
#    larger "lg" width/length variant) with no arrowhead on the tail.

$p = $ppt.ActivePresentation

# msoArrowheadNone = 1, msoArrowheadTriangle = 2
# msoArrowheadWidthMedium = 2, msoArrowheadWide = 3
# msoArrowheadLengthMedium = 2, msoArrowheadLong = 3
$msoArrowheadNone = 1
$msoArrowheadTriangle = 2
$msoArrowheadWide = 3
$msoArrowheadLong = 3

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)

    # --- 1. Rename the "ProductCustomer" cloud -> "Ownership Service" ---
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq "ProductCustomer") {
            $shp.TextFrame.TextRange.Text = "Ownership Service"
        }
    }

    # --- 2. Update connector arrowheads ---
    # The two connectors that run from "Atlas Reader" down into the two
    # cloud shapes switch from a tail-only arrow to a head-side triangle
    # (the second one using the larger "lg" width/length variant), with
    # no arrowhead on the tail.
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.Name -eq "Elbow Connector 40") {
            $shp.Line.BeginArrowheadStyle = $msoArrowheadTriangle
            $shp.Line.EndArrowheadStyle = $msoArrowheadNone
        }
        elseif ($shp.Name -eq "Elbow Connector 43") {
            $shp.Line.BeginArrowheadStyle = $msoArrowheadTriangle
            $shp.Line.BeginArrowheadWidth = $msoArrowheadWide
            $shp.Line.BeginArrowheadLength = $msoArrowheadLong
            $shp.Line.EndArrowheadStyle = $msoArrowheadNone
        }
    }
}
